# Apply tracker.xlsx edits: mark new "x" separator cells in column P (and J30/J31),
# add two new function-name rows (PageObj_GetType, Text_LoadFont) in columns O16/O17,
# add the bottom border to A28, and move the active selection to O17.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New function names appended to the O column list (rows 16 and 17)
$ws.Range("O16").Value = "PageObj_GetType"
$ws.Range("O17").Value = "Text_LoadFont"

# "x" separator marker cells in column P, using the same style as column B/D/etc (style index 2)
$markerSource = $ws.Range("B1")
$markerSource.Copy()
foreach ($r in 1..6) {
    $cell = $ws.Range("P$r")
    $cell.Value = "x"
    $cell.PasteSpecial(-4122) # xlPasteFormats
}
foreach ($r in 11..17) {
    $cell = $ws.Range("P$r")
    $cell.Value = "x"
    $cell.PasteSpecial(-4122) # xlPasteFormats
}

# New "x" separator cells in column J for rows 30-31
foreach ($r in 30..31) {
    $cell = $ws.Range("J$r")
    $cell.Value = "x"
    $cell.PasteSpecial(-4122) # xlPasteFormats
}

# A28 picks up the bottom-border style already used by the rest of row 28 (e.g. C28)
$ws.Range("C28").Copy()
$ws.Range("A28").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Selection moves to O17
$ws.Range("O17").Select()
